$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All rows refresh their "last charge end time" (column D) to the new run time.
$ws.Range("D2:D54").Value = 45969.480949074074

# Rows 19-54 get new station/terminal/last-seen data (table re-sorted after refresh).
$data = @(
    @(19, "长沙特来电飞狐四方坪东区充电站", "004A号直流", 45964.528668981482),
    @(20, "长沙特来电飞狐四方坪西区充电站", "604号直流", 45965.565891203703),
    @(21, "长沙特来电飞狐四方坪西区充电站", "603号直流", 45966.254062499997),
    @(22, "长沙特来电飞狐四方坪南区充电站", "406号直流", 45966.690613425926),
    @(23, "长沙特来电飞狐四方坪西区充电站", "504号直流", 45967.035775462966),
    @(24, "长沙特来电飞狐四方坪南区充电站", "405号直流", 45967.114155092589),
    @(25, "长沙特来电飞狐四方坪西区充电站", "505号直流", 45967.507719907408),
    @(26, "长沙特来电飞狐四方坪东区充电站", "003B号直流", 45967.530300925922),
    @(27, "长沙特来电飞狐四方坪西区充电站", "602号直流", 45967.592800925922),
    @(28, "长沙特来电飞狐四方坪南区充电站", "201号直流", 45967.63453703704),
    @(29, "长沙特来电飞狐四方坪西区充电站", "503号直流", 45968.039247685185),
    @(30, "长沙特来电飞狐四方坪西区充电站", "904号直流", 45968.065995370373),
    @(31, "长沙特来电飞狐四方坪西区充电站", "705号直流", 45968.183854166666),
    @(32, "长沙特来电飞狐四方坪西区充电站", "703号直流", 45968.23945601852),
    @(33, "长沙特来电飞狐四方坪西区充电站", "903号直流", 45968.25037037037),
    @(34, "长沙特来电飞狐四方坪西区充电站", "804号直流", 45968.268969907411),
    @(35, "长沙市开福区高岭香江国际城充电站建设项目", "106号直流", 45968.30878472222),
    @(36, "长沙特来电飞狐四方坪西区充电站", "802号直流", 45968.388923611114),
    @(37, "长沙市开福区高岭香江国际城充电站建设项目", "102号直流", 45968.454074074078),
    @(38, "长沙市开福区高岭香江国际城充电站建设项目", "103号直流", 45968.497372685182),
    @(39, "长沙市开福区高岭香江国际城充电站建设项目", "204号直流", 45968.519143518519),
    @(40, "长沙特来电飞狐四方坪东区充电站", "005A号直流", 45968.52076388889),
    @(41, "长沙特来电飞狐四方坪东区充电站", "011A号直流", 45968.521481481483),
    @(42, "长沙特来电飞狐四方坪东区充电站", "604号直流", 45968.53597222222),
    @(43, "长沙特来电飞狐四方坪西区充电站", "502号直流", 45968.550023148149),
    @(44, "长沙特来电飞狐四方坪南区充电站", "301号直流", 45968.55740740741),
    @(45, "长沙特来电飞狐四方坪南区充电站", "404号直流", 45968.557615740741),
    @(46, "长沙特来电飞狐四方坪东区充电站", "102号直流", 45968.559259259258),
    @(47, "长沙特来电飞狐四方坪南区充电站", "305号直流", 45968.563067129631),
    @(48, "长沙特来电飞狐四方坪东区充电站", "002A号直流", 45968.578900462962),
    @(49, "长沙特来电飞狐四方坪东区充电站", "905号直流", 45968.665879629632),
    @(50, "长沙市开福区高岭香江国际城充电站建设项目", "210号直流", 45968.66747685185),
    @(51, "长沙特来电飞狐四方坪西区充电站", "A05号直流", 45968.681770833333),
    @(52, "长沙市开福区高岭香江国际城充电站建设项目", "207号直流", 45968.739027777781),
    @(53, "长沙特来电飞狐四方坪南区充电站", "104号直流", 45968.753865740742),
    @(54, "长沙特来电飞狐四方坪西区充电站", "501号直流", 45968.842303240737)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$ws.Range("E6").Select()
